$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add Test ID values (column E) for each test case row
$ws.Range("E68").Value = "IT_auth_001"
$ws.Range("E69").Value = "IT_auth_002"
$ws.Range("E70").Value = "IT_auth_003"
$ws.Range("E72").Value = "IT_nav_001"
$ws.Range("E73").Value = "IT_nav_002"
$ws.Range("E75").Value = "IT_form_001"
$ws.Range("E76").Value = "IT_form_002"
$ws.Range("E77").Value = "IT_form_003"
$ws.Range("E78").Value = "IT_form_004"
$ws.Range("E80").Value = "IT_performance_001"
$ws.Range("E82").Value = "IT_layout_001"
$ws.Range("E84").Value = "IT_data_001"
$ws.Range("E85").Value = "IT_data_002"
$ws.Range("E86").Value = "IT_data_003"
$ws.Range("E88").Value = "IT_error_001"
$ws.Range("E90").Value = "IT_endpoint_001"
$ws.Range("E92").Value = "IT_get_001"
$ws.Range("E93").Value = "IT_get_002"
$ws.Range("E95").Value = "IT_psot_001"
$ws.Range("E97").Value = "IT_put_001"
$ws.Range("E99").Value = "IT_delete_001"

# Update the active cell selection to reflect where the user ended up
$ws.Range("M87").Select()
